# Auto-generated edit script for cryptos.xlsx crypto price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole affected range to Text format before writing, so that
# Excel does not auto-convert numeric-looking strings (e.g. "1.00", "3.09")
# into real numbers - the source file stores every value in B:E as text.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = "69.862.70"
$ws.Range("E2").Value = "  +1.12%  "

$ws.Range("D3").Value = "3.515.00"
$ws.Range("E3").Value = "  -0.08%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "604.63"
$ws.Range("E5").Value = "  +4.01%  "

$ws.Range("D6").Value = "170.69"
$ws.Range("E6").Value = "  -2.40%  "

$ws.Range("E7").Value = "  -0.90%  "

$ws.Range("D8").Value = "3.512.85"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").Value = "0.200"
$ws.Range("E10").Value = "  +5.62%  "

$ws.Range("D11").Value = "6.69"
$ws.Range("E11").Value = "  -0.72%  "

$ws.Range("D12").Value = "0.583"
$ws.Range("E12").Value = "  -2.97%  "

$ws.Range("D13").Value = "47.21"
$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("E14").Value = "  +0.56%  "

$ws.Range("D15").Value = "4.085.70"
$ws.Range("E15").Value = "  +0.04%  "

$ws.Range("B16").Value = "BitcoinCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D16").Value = "621.27"
$ws.Range("E16").Value = "  -7.68%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "8.38"
$ws.Range("E17").Value = "  -5.26%  "

$ws.Range("D18").Value = "3.517.63"
$ws.Range("E18").Value = "  -0.12%  "

$ws.Range("D19").Value = "69.859.91"
$ws.Range("E19").Value = "  +1.12%  "

$ws.Range("E20").Value = "  -1.91%  "

$ws.Range("D21").Value = "17.30"
$ws.Range("E21").Value = "  -1.79%  "

$ws.Range("D22").Value = "10.10"
$ws.Range("E22").Value = "  -10.49%  "

$ws.Range("D23").Value = "0.884"
$ws.Range("E23").Value = "  -2.77%  "

$ws.Range("D24").Value = "15.77"
$ws.Range("E24").Value = "  -3.26%  "

$ws.Range("D25").Value = "96.10"
$ws.Range("E25").Value = "  -2.26%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  -2.58%  "

$ws.Range("D29").Value = "9.26"
$ws.Range("E29").Value = "  -2.57%  "

$ws.Range("D30").Value = "33.13"
$ws.Range("E30").Value = "  +0.35%  "

$ws.Range("E31").Value = "  -3.94%  "

$ws.Range("D32").Value = "3.08"
$ws.Range("E32").Value = "  -4.47%  "

$ws.Range("D33").Value = "1.34"
$ws.Range("E33").Value = "  -1.87%  "

$ws.Range("D34").Value = "6.97"
$ws.Range("E34").Value = "  -4.82%  "

$ws.Range("D35").Value = "565.52"
$ws.Range("E35").Value = "  -2.44%  "

$ws.Range("D36").Value = "10.76"
$ws.Range("E36").Value = "  -1.80%  "

$ws.Range("D37").Value = "3.51"
$ws.Range("E37").Value = "  -2.50%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.101"
$ws.Range("E38").Value = "  -4.04%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "56.81"
$ws.Range("E39").Value = "  -0.94%  "

$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0453"
$ws.Range("E41").Value = "  +2.96%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.141"
$ws.Range("E42").Value = "  +3.43%  "

$ws.Range("E43").Value = "  -3.73%  "

$ws.Range("D44").Value = "3.326.63"
$ws.Range("E44").Value = "  -2.74%  "

$ws.Range("D45").Value = "33.03"
$ws.Range("E45").Value = "  -1.50%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "2.96"
$ws.Range("E46").Value = "  +2.15%  "

$ws.Range("B47").Value = "PEPE"
$ws.Range("C47").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D47").Value = "0.0₃0704"
$ws.Range("E47").Value = "  -0.68%  "

$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("E49").Value = "  -3.69%  "

$ws.Range("D50").Value = "136.04"
$ws.Range("E50").Value = "  +3.40%  "

$ws.Range("E51").Value = "  +6.08%  "

# Restore default cell style (Normal) on the edited range so the saved
# workbook does not pick up a stray Text number-format style, while the
# values themselves remain stored as text.
$editRange.Style = "Normal"
